$d = $word.ActiveDocument

# Helper: insert a new italic paragraph right after the paragraph whose
# Range currently equals $afterRange, with text $text. Mirrors what Word
# does when you place the cursor at the end of a paragraph, hit Enter,
# type text, and then apply italics to just the typed text (leaving the
# paragraph mark / pilcrow un-italicised).
function Add-ItalicParagraphAfter($afterRange, [string]$text) {
    $afterRange.InsertParagraphAfter() | Out-Null
    $newRange = $afterRange.Next(4, 1)     # wdParagraph = 4, move forward 1 paragraph
    $newRange.InsertAfter($text)
    $textRange = $d.Range($newRange.Start, $newRange.End - 1)
    $textRange.Font.Italic = $true
}

# 1) Update the activation date.
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# Work from the bottom of the document upward so earlier paragraph
# indices/ranges stay valid while we insert new paragraphs.

# 4) English translation of "Programa" (after paragraph 12).
$pPrograma = $d.Paragraphs(12).Range
Add-ItalicParagraphAfter $pPrograma "Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions."

# 3) English translation of "Programa resumido" (after paragraph 10).
$pResumido = $d.Paragraphs(10).Range
Add-ItalicParagraphAfter $pResumido "Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD)."

# 2) English translation of "Objetivos" (after paragraph 6).
$pObjetivos = $d.Paragraphs(6).Range
Add-ItalicParagraphAfter $pObjetivos "Develop knowledge in order to make the student capable of correctly interpreting the technical drawing, knowing the methodologies and tools used in the industry, giving subsidies so that they can execute, interact and modify drawings and projects throughout their professional life."
